$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update the date in A1 (45406 -> 45436) ---
$ws.Range("A1").Value = 45436

# --- Update the price list values in column D ---
$ws.Range("D33").Value = 576.08
$ws.Range("D34").Value = 630.0890000000001
$ws.Range("D35").Value = 706.599
$ws.Range("D36").Value = 859.62

# --- Re-create the merged cells so they are stored in the new order ---
$ws.Range("A1:D1").UnMerge()
$ws.Range("A9:D9").UnMerge()
$ws.Range("B33:C33").UnMerge()
$ws.Range("B34:C34").UnMerge()
$ws.Range("B36:C36").UnMerge()
$ws.Range("B32:C32").UnMerge()
$ws.Range("A11:D11").UnMerge()
$ws.Range("A10:D10").UnMerge()
$ws.Range("B35:C35").UnMerge()

$ws.Range("A1:D1").Merge()
$ws.Range("A9:D9").Merge()
$ws.Range("B33:C33").Merge()
$ws.Range("B34:C34").Merge()
$ws.Range("B36:C36").Merge()
$ws.Range("B32:C32").Merge()
$ws.Range("A11:D11").Merge()
$ws.Range("A10:D10").Merge()
$ws.Range("B35:C35").Merge()
